# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Map of row -> new F value for sheet "展览" (rows 3-40)
$updates1 = @{
    3  = 563
    4  = 1135
    5  = 122
    6  = 79
    9  = 1168
    10 = 16438
    12 = 204
    13 = 1040
    14 = 6385
    15 = 640
    17 = 82
    18 = 28
    21 = 44
    24 = 39
    25 = 23
    28 = 226
    29 = 898
    30 = 58
    31 = 5060
    32 = 504
    33 = 11362
    36 = 153
    37 = 212
    38 = 3843
    39 = 271
    40 = 73
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

# Map of row -> new F value for sheet "全部类型" (rows 3-41, has an extra row at 33)
$updates4 = @{
    3  = 563
    4  = 1135
    5  = 122
    6  = 79
    9  = 1168
    10 = 16438
    12 = 204
    13 = 1040
    14 = 6385
    15 = 640
    17 = 82
    18 = 28
    21 = 44
    24 = 39
    25 = 23
    28 = 226
    29 = 898
    30 = 58
    31 = 5060
    32 = 504
    34 = 11362
    37 = 153
    38 = 212
    39 = 3843
    40 = 271
    41 = 73
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
